# Apply the "TS Parameters" sheet updates:
#  - G1 header: "Mapping Result" -> "Mapping Results"
#  - Most "None" placeholders in column G -> a single space " "
#  - G49 "None" -> a specific JSONata error message

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TS Parameters")

# Update header text
$ws.Range("G1").Value = "Mapping Results"

# Cells whose "None" placeholder becomes a single space
$spaceCells = @(
    "G8","G9","G10","G11","G12","G13","G14","G15","G16","G17","G18","G19","G20",
    "G24","G25",
    "G27","G28",
    "G32","G33","G34","G35","G36","G37","G38","G39","G40",
    "G42","G43","G44","G45","G46","G47","G48",
    "G50","G51","G52",
    "G54","G55","G56",
    "G59","G60"
)

foreach ($cellRef in $spaceCells) {
    $ws.Range($cellRef).Value = " "
}

# Special cell: error message for JSONata expression
$ws.Range("G49").Value = 'Error in expression for Clinical Study Sponsor; Sponsor; Study Sponsor: study.versions.($sponsorIdVal:=roles[code.code="C70793"].organizationIds[0];'
